# "date formatting for export"
#
# The template's {lastmod} placeholder becomes {lastmod:date}, and a new
# "ref. version" / {refmod:date} row is added right below it so the export
# can carry a second (reference) last-modified date, each rendered with
# explicit date formatting. The legend in the top-right corner (added /
# changed / linked / unlinked swatches) gets its "linked"/"unlinked" text
# styles swapped: the former purple italic-underline becomes a plain blue
# underline (hyperlink-style), and the former plain grey becomes an italic
# grey. Finally the frozen-pane scroll position / selection is nudged to
# rest on the legend.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D3: {lastmod} -> {lastmod:date}, now styled to match the other
#     italic placeholder cells (C4/D4) ---
$ws.Range("D3").Value = "{lastmod:date}"
$ws.Range("D3").Font.Italic = $true

# --- C4/D4: new "ref. version" / "{refmod:date}" row ---
$ws.Range("C4").Value = "ref. version"
$ws.Range("C4").WrapText = $false
$ws.Range("C4").HorizontalAlignment = -4152
$ws.Range("C4").Font.Italic = $true

$ws.Range("D4").Value = "{refmod:date}"

# --- Legend swatches, column O/P (row2 = added/linked, row3 = changed/unlinked) ---
# P3 ("unlinked" sample) becomes italic, keeping its muted grey theme color.
$ws.Range("P3").Font.Italic = $true

# P2 ("linked" sample) drops the italic and swaps its purple underline for a
# plain blue (FF0070C0) underline - a hyperlink look.
$ws.Range("P2").Font.Italic = $false
$ws.Range("P2").Font.Color = 12611584

# --- View state: select the legend swatches so the frozen pane scrolls to
#     show row 7 at the top again ---
$ws.Range("P2:P3").Select()
